$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "ShowName" column (U) with header and per-row placeholder values,
# and strip the ".png" extension that was mistakenly baked into the Icon (T) column strings.
$ws.Range("U1").Value = "ShowName"
$ws.Range("T2").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U2").Value = "ShowName_1"
$ws.Range("T3").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U3").Value = "ShowName_2"
$ws.Range("T4").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U4").Value = "ShowName_3"
$ws.Range("T5").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U5").Value = "ShowName_4"
$ws.Range("T6").Value = "UI/SteampunkUI/resource/icons/icon_Coin"
$ws.Range("U6").Value = "ShowName_5"
$ws.Range("T7").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U7").Value = "ShowName_6"
$ws.Range("T8").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U8").Value = "ShowName_7"
$ws.Range("T9").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U9").Value = "ShowName_8"
$ws.Range("T10").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U10").Value = "ShowName_9"
$ws.Range("T11").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U11").Value = "ShowName_10"
$ws.Range("T12").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U12").Value = "ShowName_11"
$ws.Range("T13").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U13").Value = "ShowName_12"
$ws.Range("T14").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U14").Value = "ShowName_13"
$ws.Range("T15").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U15").Value = "ShowName_14"
$ws.Range("T16").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U16").Value = "ShowName_15"
$ws.Range("T17").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U17").Value = "ShowName_16"
$ws.Range("T18").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U18").Value = "ShowName_17"
$ws.Range("T19").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U19").Value = "ShowName_18"
$ws.Range("T20").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U20").Value = "ShowName_19"
$ws.Range("T21").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U21").Value = "ShowName_20"
$ws.Range("T22").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U22").Value = "ShowName_21"
$ws.Range("T23").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U23").Value = "ShowName_22"
$ws.Range("T24").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U24").Value = "ShowName_23"
$ws.Range("T25").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U25").Value = "ShowName_24"
$ws.Range("T26").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U26").Value = "ShowName_25"
$ws.Range("T27").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U27").Value = "ShowName_26"
$ws.Range("T28").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U28").Value = "ShowName_27"
$ws.Range("T29").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U29").Value = "ShowName_28"
$ws.Range("T30").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U30").Value = "ShowName_29"
$ws.Range("T31").Value = "UI/SteampunkUI/resource/icons/Icon_Gold_few"
$ws.Range("U31").Value = "ShowName_30"
$ws.Range("T32").Value = "UI/SteampunkUI/resource/icons/Icon_Gold_middle"
$ws.Range("U32").Value = "ShowName_31"
$ws.Range("T33").Value = "UI/SteampunkUI/resource/icons/Icon_Gold_lot"
$ws.Range("U33").Value = "ShowName_32"
$ws.Range("T34").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U34").Value = "ShowName_33"
$ws.Range("T35").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U35").Value = "ShowName_34"
$ws.Range("T36").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U36").Value = "ShowName_35"
$ws.Range("T37").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U37").Value = "ShowName_36"
$ws.Range("T38").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U38").Value = "ShowName_37"
$ws.Range("T39").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U39").Value = "ShowName_38"
$ws.Range("T40").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U40").Value = "ShowName_39"
$ws.Range("T41").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U41").Value = "ShowName_40"
$ws.Range("T42").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U42").Value = "ShowName_41"
$ws.Range("T43").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U43").Value = "ShowName_42"
$ws.Range("T44").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U44").Value = "ShowName_43"
$ws.Range("T45").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U45").Value = "ShowName_44"
$ws.Range("T46").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U46").Value = "ShowName_45"
$ws.Range("T47").Value = "UI/SteampunkUI/resource/icons/img_equip"
$ws.Range("U47").Value = "ShowName_46"

# Restore the view: keep the existing freeze at row 1 / column A, scroll the
# frozen pane so column P is visible, and leave the cursor on the last data row.
$win = $excel.ActiveWindow
$win.ScrollColumn = 16
$win.ScrollRow = 2
$ws.Range("W43").Select()
